# quarterly.xlsx update: drop the oldest quarter column (D) and append the
# newest quarter (new column M) - "update database and change read_price algorithm"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Drop the oldest quarter (column D), shifting every later quarter one
#    column to the left (D<-E, E<-F, ... L<-M).
$ws.Columns.Item(4).Delete()

# 2) Set the column width for the freshly appended column M (matches the
#    width of the other "publish date" columns E/I -- raw XML width 31,
#    which the ColumnWidth COM property reports/accepts as 30.17).
$ws.Columns.Item(13).ColumnWidth = 30.17

# 3) Header row 8: new quarter label.
$ws.Range("M8").Value = "فصل چهارم منتهی به 1401/12"

# 4) Header row 9: new quarter's publish date.
$ws.Range("M9").Value = "1402-02-27"

# 5) New column M data for every data row (the latest quarter's figures).
$ws.Range("M12").Value = -6350150
$ws.Range("M13").Value = 326200
$ws.Range("M14").Value = -6023950

$ws.Range("M16").Value = 0
$ws.Range("M17").Value = -963434
$ws.Range("M18").Value = 0
$ws.Range("M19").Value = 0
$ws.Range("M20").Value = 0
$ws.Range("M21").Value = 0
$ws.Range("M22").Value = 0
$ws.Range("M23").Value = 0
$ws.Range("M24").Value = 0
$ws.Range("M25").Value = 0
$ws.Range("M26").Value = -1220000
$ws.Range("M27").Value = 0
$ws.Range("M28").Value = 0
$ws.Range("M29").Value = 0
$ws.Range("M30").Value = 0
$ws.Range("M31").Value = 244259
$ws.Range("M32").Value = -1939175
$ws.Range("M33").Value = -7963125

$ws.Range("M35").Value = 0
# Row 36 ("cash received from share premium"): the two periods with actual
# (zero) data now show a literal 0 instead of the "-" placeholder.
$ws.Range("I36").Value = 0
$ws.Range("M36").Value = 0

$ws.Range("M37").Value = 0
$ws.Range("M38").Value = 0
$ws.Range("M39").Value = 16099860
$ws.Range("M40").Value = -14039582
$ws.Range("M41").Value = -413962
$ws.Range("M42").Value = 0
$ws.Range("M43").Value = 0
$ws.Range("M44").Value = 0
$ws.Range("M45").Value = 0
$ws.Range("M46").Value = 0
$ws.Range("M47").Value = 0
$ws.Range("M48").Value = 0
$ws.Range("M49").Value = 0
$ws.Range("M50").Value = 0
$ws.Range("M51").Value = 1646316
$ws.Range("M52").Value = -6316809
$ws.Range("M53").Value = 7994029
$ws.Range("M54").Value = -12454
$ws.Range("M55").Value = 1664766
